# Auto-generated edit script to apply metrics corrections
# "fix - metricas de bugs e dt de cada tribo agora estao sendo somadas corretamente"

$wb = $excel.ActiveWorkbook

$wsOctubre = $wb.Worksheets.Item("Octubre")
$wsNoviembre = $wb.Worksheets.Item("Noviembre")

# --- Octubre sheet updates ---
$wsOctubre.Range("I4").Value = 0
$wsOctubre.Range("J6").Value = 19
$wsOctubre.Range("I7").Value = 31
$wsOctubre.Range("I10").Value = 0
$wsOctubre.Range("J10").Value = 0
$wsOctubre.Range("I11").Value = 46
$wsOctubre.Range("I12").Value = 22
$wsOctubre.Range("J12").Value = 19

# --- Noviembre sheet updates ---
$wsNoviembre.Range("B2").Value = 3131
$wsNoviembre.Range("C2").Value = 1887
$wsNoviembre.Range("D2").Value = 214
$wsNoviembre.Range("E2").Value = 1030
$wsNoviembre.Range("G2").Value = 3.67
$wsNoviembre.Range("H2").Value = 15.46
$wsNoviembre.Range("B3").Value = 13344
$wsNoviembre.Range("C3").Value = 944
$wsNoviembre.Range("D3").Value = 343
$wsNoviembre.Range("E3").Value = 829
$wsNoviembre.Range("F3").Value = 11228
$wsNoviembre.Range("G3").Value = 9.82
$wsNoviembre.Range("H3").Value = 22.08
$wsNoviembre.Range("I3").Value = 1
$wsNoviembre.Range("J3").Value = 2
$wsNoviembre.Range("B4").Value = 13812
$wsNoviembre.Range("C4").Value = 606
$wsNoviembre.Range("F4").Value = 12702
$wsNoviembre.Range("G4").Value = 16.63
$wsNoviembre.Range("H4").Value = 6.2
$wsNoviembre.Range("I4").Value = 0
$wsNoviembre.Range("I5").Value = 0
$wsNoviembre.Range("B6").Value = 9397
$wsNoviembre.Range("C6").Value = 1769
$wsNoviembre.Range("D6").Value = 275
$wsNoviembre.Range("E6").Value = 918
$wsNoviembre.Range("G6").Value = 6.89
$wsNoviembre.Range("H6").Value = 5.98
$wsNoviembre.Range("I6").Value = 8
$wsNoviembre.Range("J6").Value = 21
$wsNoviembre.Range("B7").Value = 24502
$wsNoviembre.Range("C7").Value = 2752
$wsNoviembre.Range("D7").Value = 759
$wsNoviembre.Range("E7").Value = 3549
$wsNoviembre.Range("G7").Value = 15.68
$wsNoviembre.Range("H7").Value = 6.17
$wsNoviembre.Range("J7").Value = 21
$wsNoviembre.Range("B8").Value = 12430
$wsNoviembre.Range("C8").Value = 938
$wsNoviembre.Range("D8").Value = 564
$wsNoviembre.Range("E8").Value = 1917
$wsNoviembre.Range("G8").Value = 10.56
$wsNoviembre.Range("H8").Value = 4.16
$wsNoviembre.Range("I8").Value = 16
$wsNoviembre.Range("B9").Value = 11024
$wsNoviembre.Range("C9").Value = 3795
$wsNoviembre.Range("D9").Value = 559
$wsNoviembre.Range("E9").Value = 1222
$wsNoviembre.Range("G9").Value = 9.369999999999999
$wsNoviembre.Range("H9").Value = 6.69
$wsNoviembre.Range("B10").Value = 22161
$wsNoviembre.Range("C10").Value = 328
$wsNoviembre.Range("D10").Value = 56
$wsNoviembre.Range("E10").Value = 270
$wsNoviembre.Range("F10").Value = 21507
$wsNoviembre.Range("G10").Value = 5.6
$wsNoviembre.Range("H10").Value = 19.01
$wsNoviembre.Range("I10").Value = 0
$wsNoviembre.Range("J10").Value = 0
$wsNoviembre.Range("B11").Value = 16914
$wsNoviembre.Range("C11").Value = 1432
$wsNoviembre.Range("D11").Value = 532
$wsNoviembre.Range("E11").Value = 1899
$wsNoviembre.Range("F11").Value = 13051
$wsNoviembre.Range("G11").Value = 14.29
$wsNoviembre.Range("H11").Value = 5.51
$wsNoviembre.Range("J11").Value = 22
$wsNoviembre.Range("B12").Value = 17055
$wsNoviembre.Range("C12").Value = 1270
$wsNoviembre.Range("D12").Value = 127
$wsNoviembre.Range("E12").Value = 702
$wsNoviembre.Range("F12").Value = 14955
$wsNoviembre.Range("G12").Value = 9.359999999999999
$wsNoviembre.Range("I12").Value = 8
$wsNoviembre.Range("J12").Value = 21
$wsNoviembre.Range("B13").Value = 7077
$wsNoviembre.Range("C13").Value = 2841
$wsNoviembre.Range("D13").Value = 386
$wsNoviembre.Range("E13").Value = 1126
$wsNoviembre.Range("G13").Value = 6.52
$wsNoviembre.Range("H13").Value = 11.08
$wsNoviembre.Range("B14").Value = 13344
$wsNoviembre.Range("C14").Value = 944
$wsNoviembre.Range("D14").Value = 343
$wsNoviembre.Range("E14").Value = 829
$wsNoviembre.Range("F14").Value = 11228
$wsNoviembre.Range("G14").Value = 9.82
$wsNoviembre.Range("H14").Value = 22.08
$wsNoviembre.Range("I14").Value = 1
$wsNoviembre.Range("J14").Value = 2
